# Refresh the scraped cryptocurrency price snapshot (and a couple of the
# derived "WorstIn24h"/ticker labels in column E) to match the latest
# GitHub Actions run. Column D values are numeric-looking text (price
# quotes captured as strings), so they are written with a leading
# apostrophe (quote-prefix) to force Excel to store them as text and keep
# exact formatting (e.g. trailing zeros such as "3.850" or "0.1380").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'243.62"
$ws.Range("D3").Formula = "'23.71"
$ws.Range("D4").Formula = "'5.257"
$ws.Range("D5").Formula = "'0.05805"
$ws.Range("D6").Formula = "'6.485"
$ws.Range("D8").Formula = "'0.8083"
$ws.Range("D9").Formula = "'0.8723"
$ws.Range("D10").Formula = "'0.1380"
$ws.Range("D11").Formula = "'0.07273"
$ws.Range("D12").Formula = "'0.03065"
$ws.Range("D13").Formula = "'0.03052"
$ws.Range("D14").Formula = "'0.09309"
$ws.Range("D15").Formula = "'3.850"
$ws.Range("D16").Formula = "'0.001531"
$ws.Range("D17").Formula = "'0.04694"
$ws.Range("D18").Formula = "'0.0006049"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("D19").Formula = "'0.006179"
$ws.Range("D21").Formula = "'0.004590"
$ws.Range("D22").Formula = "'0.00008699"
$ws.Range("D25").Formula = "'0.3208"
$ws.Range("D28").Formula = "'0.0002344"
$ws.Range("D40").Formula = "'0.03790"
$ws.Range("D41").Formula = "'0.006349"
$ws.Range("D43").Formula = "'0.002600"
$ws.Range("D44").Formula = "'0.006858"
$ws.Range("E44").Value = "43LocalTradersLCTWorstin24h"
$ws.Range("D45").Formula = "'0.00005485"
$ws.Range("D47").Formula = "'0.5499"
$ws.Range("D48").Formula = "'0.006828"
